# Updates the date header and all multiplication problems/answers in the
# table to a newer set of generated values.
$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-02-20 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-02-21 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("445×4=1780", $true, $false, $false, $false, $false, $true, 1, $false, "108×8=864", 2) | Out-Null
$d.Content.Find.Execute("395×6=2370", $true, $false, $false, $false, $false, $true, 1, $false, "809×5=4045", 2) | Out-Null
$d.Content.Find.Execute("239×3=717", $true, $false, $false, $false, $false, $true, 1, $false, "378×7=2646", 2) | Out-Null
$d.Content.Find.Execute("942×5=4710", $true, $false, $false, $false, $false, $true, 1, $false, "530×7=3710", 2) | Out-Null
$d.Content.Find.Execute("433×7=3031", $true, $false, $false, $false, $false, $true, 1, $false, "299×7=2093", 2) | Out-Null
$d.Content.Find.Execute("203×6=1218", $true, $false, $false, $false, $false, $true, 1, $false, "421×4=1684", 2) | Out-Null
$d.Content.Find.Execute("681×8=5448", $true, $false, $false, $false, $false, $true, 1, $false, "225×6=1350", 2) | Out-Null
$d.Content.Find.Execute("248×7=1736", $true, $false, $false, $false, $false, $true, 1, $false, "322×5=1610", 2) | Out-Null
$d.Content.Find.Execute("998×7=6986", $true, $false, $false, $false, $false, $true, 1, $false, "889×6=5334", 2) | Out-Null
$d.Content.Find.Execute("802×6=4812", $true, $false, $false, $false, $false, $true, 1, $false, "837×9=7533", 2) | Out-Null
$d.Content.Find.Execute("776×9=6984", $true, $false, $false, $false, $false, $true, 1, $false, "635×6=3810", 2) | Out-Null
$d.Content.Find.Execute("234×6=1404", $true, $false, $false, $false, $false, $true, 1, $false, "336×2=672", 2) | Out-Null
$d.Content.Find.Execute("177×4=708", $true, $false, $false, $false, $false, $true, 1, $false, "221×7=1547", 2) | Out-Null
$d.Content.Find.Execute("796×5=3980", $true, $false, $false, $false, $false, $true, 1, $false, "958×9=8622", 2) | Out-Null
$d.Content.Find.Execute("990×5=4950", $true, $false, $false, $false, $false, $true, 1, $false, "227×8=1816", 2) | Out-Null
$d.Content.Find.Execute("249×5=1245", $true, $false, $false, $false, $false, $true, 1, $false, "689×8=5512", 2) | Out-Null
$d.Content.Find.Execute("120×8=960", $true, $false, $false, $false, $false, $true, 1, $false, "673×6=4038", 2) | Out-Null
$d.Content.Find.Execute("176×7=1232", $true, $false, $false, $false, $false, $true, 1, $false, "522×4=2088", 2) | Out-Null
$d.Content.Find.Execute("867×8=6936", $true, $false, $false, $false, $false, $true, 1, $false, "116×8=928", 2) | Out-Null
$d.Content.Find.Execute("331×4=1324", $true, $false, $false, $false, $false, $true, 1, $false, "345×2=690", 2) | Out-Null
$d.Content.Find.Execute("851×9=7659", $true, $false, $false, $false, $false, $true, 1, $false, "537×7=3759", 2) | Out-Null
$d.Content.Find.Execute("604×7=4228", $true, $false, $false, $false, $false, $true, 1, $false, "113×4=452", 2) | Out-Null
$d.Content.Find.Execute("224×4=896", $true, $false, $false, $false, $false, $true, 1, $false, "112×2=224", 2) | Out-Null
$d.Content.Find.Execute("494×4=1976", $true, $false, $false, $false, $false, $true, 1, $false, "679×3=2037", 2) | Out-Null
$d.Content.Find.Execute("460×4=1840", $true, $false, $false, $false, $false, $true, 1, $false, "849×6=5094", 2) | Out-Null
